# Fixed a bug in ReplaceSymbol
# The rows A5:F20 hold (symbol, count, c, d, e, f) tuples that were associated
# with the wrong symbol values. This corrects the row ordering so that each
# row's counts line up with the correct symbol in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 5
$lastRow = 20

# Capture the current contents of every row (columns A:F) before mutating
# anything, since several rows trade places with each other.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 6; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Maps destination row -> source row (i.e. destination row should end up
# holding what is currently in the source row).
$mapping = @{
    5  = 8
    6  = 7
    7  = 5
    8  = 10
    9  = 15
    10 = 9
    11 = 12
    12 = 13
    13 = 11
    14 = 6
    15 = 14
    16 = 20
    17 = 17
    18 = 19
    19 = 18
    20 = 16
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
